# Update the USD Amount ("T2") figure on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("T2").Value = 234600
